$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Columns.Item(9).Insert()
$ws.Range("I1").Value = "Count"
$ws.Range("I2").Value = "string"
$ws.Range("I3").Value = "FALSE"
$ws.Range("I4").Value = "FALSE"
$ws.Range("I5").Value = "TRUE"
$ws.Range("I6").Value = "FALSE"
$ws.Range("I7").Value = "FALSE"
$ws.Range("I8").Value = "FALSE"

for ($r = 10; $r -le 71; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
}

$ws.Range("I9").ClearContents()
